# Regenerate orders with updated distance/size codes.
#
# The stimulus-naming scheme embeds a Distance code (D64/D51/D80) and a
# Size code (S30, alongside unaffected S20/S25) inside several columns:
#   Condition, Filename_Left, Filename_Right, Distance, Size
# This run renumbers those codes:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# Every occurrence (whole-value columns like Distance/Size as well as the
# substrings embedded in Condition/Filename_* values) needs to change, so a
# simple ordered string replace on each affected cell handles every case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# Column letters (and numbers) that can contain the distance/size codes:
#   B = Condition, D = Filename_Left, E = Filename_Right,
#   H = Distance, J = Size
$affectedColumns = @(2, 4, 5, 8, 10)

$changed = 0
for ($row = 2; $row -le $lastRow; $row++) {
    foreach ($col in $affectedColumns) {
        $cell = $ws.Cells.Item($row, $col)
        $current = $cell.Value2
        if ($current -ne $null) {
            $text = [string]$current
            $updated = $text.Replace("D64", "D69").Replace("D51", "D55").Replace("D80", "D86").Replace("S30", "S31")
            if ($updated -ne $text) {
                $cell.Value = $updated
                $changed++
            }
        }
    }
}

Write-Output "Updated $changed cells with new distance/size codes."
